$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.217.03"
$ws.Range("E2").Value = "  +2.43%  "

$ws.Range("D3").Value = "1.998.47"
$ws.Range("E3").Value = "  +6.53%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "'0.7771"
$ws.Range("E5").Value = "  +64.87%  "

$ws.Range("D6").Value = "'256.74"
$ws.Range("E6").Value = "  +4.35%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "'0.3466"
$ws.Range("E8").Value = "  +20.69%  "

$ws.Range("D9").Value = "'28.71"
$ws.Range("E9").Value = "  +31.03%  "

$ws.Range("D10").Value = "'0.07004"
$ws.Range("E10").Value = "  +7.56%  "

$ws.Range("D11").Value = "'0.8536"
$ws.Range("E11").Value = "  +17.26%  "

$ws.Range("D12").Value = "'0.08202"
$ws.Range("E12").Value = "  +5.07%  "

$ws.Range("D13").Value = "1.997.23"
$ws.Range("E13").Value = "  +6.49%  "

$ws.Range("D14").Value = "'100.69"
$ws.Range("E14").Value = "  -0.09%  "

$ws.Range("D15").Value = "'5.598"
$ws.Range("E15").Value = "  +8.16%  "

$ws.Range("D16").Value = "'15.56"
$ws.Range("E16").Value = "  +18.85%  "

$ws.Range("D17").Value = "'274.11"
$ws.Range("E17").Value = "  -3.75%  "

$ws.Range("D18").Value = "31.225.82"
$ws.Range("E18").Value = "  +2.46%  "

$ws.Range("D19").Value = "'5.957"
$ws.Range("E19").Value = "  +11.58%  "

$ws.Range("D20").Value = "'0.000007909"
$ws.Range("E20").Value = "  +5.63%  "

$ws.Range("D21").Value = "2.260.19"
$ws.Range("E21").Value = "  +6.82%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").Value = "'7.109"
$ws.Range("E24").Value = "  +12.01%  "

$ws.Range("D25").Value = "'10.03"
$ws.Range("E25").Value = "  +10.78%  "

$ws.Range("D26").Value = "'164.75"
$ws.Range("E26").Value = "  +1.65%  "

$ws.Range("D27").Value = "'0.1468"
$ws.Range("E27").Value = "  +51.64%  "

$ws.Range("D28").Value = "'19.92"
$ws.Range("E28").Value = "  +4.89%  "

$ws.Range("D29").Value = "'2.325"
$ws.Range("E29").Value = "  +22.38%  "

$ws.Range("D30").Value = "'1.607"
$ws.Range("E30").Value = "  +7.27%  "

$ws.Range("D31").Value = "'4.611"
$ws.Range("E31").Value = "  +9.04%  "

$ws.Range("D32").Value = "'1.356"
$ws.Range("E32").Value = "  +2.59%  "

$ws.Range("D33").Value = "'4.429"
$ws.Range("E33").Value = "  +6.69%  "

$ws.Range("D34").Value = "'0.05211"
$ws.Range("E34").Value = "  +8.34%  "

$ws.Range("E35").Value = "  +9.12%  "

$ws.Range("D36").Value = "'0.7739"
$ws.Range("E36").Value = "  +12.04%  "

$ws.Range("D37").Value = "'2.753"
$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("D38").Value = "'0.02001"
$ws.Range("E38").Value = "  +5.13%  "

$ws.Range("E39").Value = "  +2.60%  "

$ws.Range("D40").Value = "'6.711"
$ws.Range("E40").Value = "  +6.41%  "

$ws.Range("D41").Value = "'79.43"
$ws.Range("E41").Value = "  +4.75%  "

$ws.Range("D42").Value = "'2.145"
$ws.Range("E42").Value = "  +9.31%  "

$ws.Range("D43").Value = "'0.4683"
$ws.Range("E43").Value = "  +10.90%  "

$ws.Range("D44").Value = "'106.03"
$ws.Range("E44").Value = "  +5.01%  "

$ws.Range("E45").Value = "  +2.73%  "

$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.23%  "

$ws.Range("D47").Value = "'7.718"
$ws.Range("E47").Value = "  +10.03%  "

$ws.Range("D48").Value = "'9.902"
$ws.Range("E48").Value = "  +1.69%  "

$ws.Range("D51").Value = "'1.526"
$ws.Range("E51").Value = "  +14.74%  "

# Row 49/50: swap Decentraland/Elrond coin entries with updated prices
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'36.81"
$ws.Range("E49").Value = "  +4.91%  "

$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").Value = "'0.4307"
$ws.Range("E50").Value = "  +9.80%  "

